$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings / cell values for the new "Contato" table (column D) ---
$ws.Range("D1").Value = "Contato"
$ws.Range("D2").Value = "Id_Contato"
$ws.Range("D3").Value = "Nome"
$ws.Range("D4").Value = "Email"
$ws.Range("D5").Value = "Empresa"
$ws.Range("D6").Value = "Problema"

# --- Formatting: reuse existing styles where possible by copying formats ---

# D1 header cell: same style as the other header cells (dark fill, border, centered)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# D2:D6 start from the plain thin-bordered style already present on C4/C5 (D2 gets a
# fill colour on top of it below; D3:D6 stay exactly like this)
$ws.Range("C4").Copy() | Out-Null
$ws.Range("D2:D6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# C2:C5 already carry the right fills (yellow/green/none) - they just lose their right
# border now that column D sits right beside them
$ws.Range("C2:C5").Borders.Item(10).LineStyle = -4142   # xlEdgeRight / xlLineStyleNone

# D2 additionally gets a brand new highlight fill (theme "Orange, Accent 2, Lighter 40%")
$d2 = $ws.Range("D2")
$d2.Interior.ThemeColor = 6
$d2.Interior.TintAndShade = 0.39997558519241921

Write-Host "done"
